$d = $word.ActiveDocument

# Helper pattern: replace all exact occurrences of a line-number reference.
# Replace=2 means wdReplaceAll (safe because each search string below is
# crafted to be unique/unambiguous for its intended occurrence(s)).

# 1) caseConditional(M2DocEvaluator.java:1438) -> 1477  (first occurrence only)
$d.Content.Find.Execute(
    "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1438)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1477)", 2) | Out-Null

# 2) doSwitch(M2DocEvaluator.java:1216) -> 1239  (all 5 occurrences)
$d.Content.Find.Execute(
    "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)", 2) | Out-Null

# 3) caseBlock(M2DocEvaluator.java:1425) -> 1464  (both occurrences)
$d.Content.Find.Execute(
    "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)", 2) | Out-Null

# 4) caseConditional(M2DocEvaluator.java:1449) -> 1488  (second occurrence only)
$d.Content.Find.Execute(
    "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1449)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1488)", 2) | Out-Null

# 5) caseDocumentTemplate(M2DocEvaluator.java:287) -> 296
$d.Content.Find.Execute(
    "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)", 2) | Out-Null

# 6) generate(M2DocEvaluator.java:276) -> 281
$d.Content.Find.Execute(
    "M2DocEvaluator.generate(M2DocEvaluator.java:276)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocEvaluator.generate(M2DocEvaluator.java:281)", 2) | Out-Null

# 7) M2DocUtils.generate(M2DocUtils.java:694) -> 805
$d.Content.Find.Execute(
    "M2DocUtils.generate(M2DocUtils.java:694)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocUtils.generate(M2DocUtils.java:805)", 2) | Out-Null

# 8) prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480) -> 511
$d.Content.Find.Execute(
    "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)", $true, $false, $false, $false, $false,
    $true, 1, $false, "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)", 2) | Out-Null

# 9) generation(AbstractTemplatesTestSuite.java:389) -> 420
$d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)", $true, $false, $false, $false, $false,
    $true, 1, $false, "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)", 2) | Out-Null

# 10) Insert a new stack frame line for RunBefores before the second RunAfters
#     occurrence (the one directly following "ParentRunner.java:268)").
$search = "ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
$replace = "ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
$d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 1) | Out-Null

Write-Host "Done"
